# fall 22 week 9 complete
# Appends 23 new Player_1/Points_1/Player_2/Points_2 matchup rows to Sheet1,
# immediately after the existing data (which ends at row 1361).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(3,0,7,3),
    @(4,2,3,1),
    @(4,2,5,1),
    @(3,2,3,1),
    @(6,3,2,0),
    @(3,0,2,3),
    @(6,2,5,0),
    @(4,0,3,3),
    @(7,1,6,2),
    @(3,1,4,2),
    @(5,2,5,0),
    @(3,2,4,1),
    @(4,0,6,3),
    @(4,2,2,1),
    @(5,0,6,2),
    @(4,0,3,3),
    @(6,0,6,2),
    @(4,0,3,3),
    @(4,2,2,1),
    @(5,3,3,0),
    @(5,2,3,1),
    @(4,0,5,3),
    @(6,0,7,2)
)

$startRow = 1362
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$lastRow = $r - 1
$nextRow = $r

# Scroll/position the view the way Excel would after typing this block of
# rows in (mirrors the workbook's saved sheetView: new active cell one row
# past the last entry, with the view scrolled down to keep it in frame).
$excel.ActiveWindow.ScrollRow = ($startRow - 1)
$ws.Range("A" + $nextRow).Select()
